$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (AD1:AF1), copying the header style (bold, bordered,
# centered) from an existing header cell so the new columns match the rest
# of the header row.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every data row.
for ($r = 2; $r -le 43; $r++) {
    $ws.Range("AD$r").Value = 86
    $ws.Range("AE$r").Value = 75
    $ws.Range("AF$r").Value = 0
}
